$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily format as text so the numeric-looking strings "7" and "2"
# are stored as shared-string text (matching the source data), not
# auto-converted to numbers. Clear the format again afterwards so no
# stray cell style is left behind.
$ws.Range("B1:B2").NumberFormat = "@"
$ws.Range("B1").Value = "7"
$ws.Range("B2").Value = "2"
$ws.Range("B1:B2").ClearFormats()

$ws.Range("N5").Select() | Out-Null
